$d = $word.ActiveDocument

# Update the date line (first paragraph)
$d.Paragraphs.Item(1).Range.Text = "2023-03-24 Friday"

# Update each answer cell in the table (row-major order, 20 rows x 5 cols)
$t = $d.Tables.Item(1)
$answers = @(
    "56-24=32",
    "18+76=94",
    "24+12=36",
    "87-80=7",
    "80-61=19",
    "77-7=70",
    "39+1=40",
    "14+52=66",
    "54-12=42",
    "51-43=8",
    "58+36=94",
    "33+34=67",
    "17+4=21",
    "18+24=42",
    "53-46=7",
    "40+59=99",
    "33+54=87",
    "63+35=98",
    "35+44=79",
    "14+0=14",
    "71-44=27",
    "82-53=29",
    "2+24=26",
    "42+48=90",
    "39-18=21",
    "92-12=80",
    "78-55=23",
    "82+10=92",
    "20+5=25",
    "17+19=36",
    "91-0=91",
    "79-35=44",
    "50-33=17",
    "31-4=27",
    "28+71=99",
    "88-38=50",
    "87+10=97",
    "62-50=12",
    "6+20=26",
    "68+20=88",
    "96-5=91",
    "79-31=48",
    "26-24=2",
    "94-36=58",
    "4+1=5",
    "41+38=79",
    "16+39=55",
    "1+38=39",
    "86-43=43",
    "67+13=80",
    "65-27=38",
    "62+35=97",
    "25+69=94",
    "10+77=87",
    "31+17=48",
    "61-58=3",
    "2+6=8",
    "47-15=32",
    "3+88=91",
    "35+60=95",
    "55+22=77",
    "79-1=78",
    "21+58=79",
    "64-31=33",
    "32+20=52",
    "52-13=39",
    "8+31=39",
    "49-49=0",
    "52+24=76",
    "28-4=24",
    "22+1=23",
    "38-23=15",
    "94-38=56",
    "41+7=48",
    "86-66=20",
    "84+4=88",
    "68+2=70",
    "64-54=10",
    "67-67=0",
    "64-4=60",
    "31+14=45",
    "49+34=83",
    "74-63=11",
    "61+21=82",
    "90+3=93",
    "62+30=92",
    "22-10=12",
    "52-38=14",
    "10+54=64",
    "70-9=61",
    "75-62=13",
    "18+71=89",
    "8+45=53",
    "1+32=33",
    "59+13=72",
    "64-3=61",
    "92-44=48",
    "94-21=73",
    "37-36=1",
    "51+22=73"
)

$idx = 0
for ($r = 1; $r -le 20; $r++) {
    for ($c = 1; $c -le 5; $c++) {
        $t.Cell($r, $c).Range.Text = $answers[$idx]
        $idx = $idx + 1
    }
}

Write-Host "Done. Updated $idx cells."